$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update group N2 (row 3) score values
$ws.Range("B3").Value = "400/400"
$ws.Range("D3").Value = 8

# Update group N3 (row 4) score values - becomes textual "400/400" instead of numeric 0
# Leading apostrophe preserves the existing quote-prefix cell style (s="1")
$ws.Range("B4").Value = "'400/400"
$ws.Range("D4").Value = 8

# Update group N5 (row 6) score values
$ws.Range("B6").Value = "400/400"
$ws.Range("D6").Value = 9

# Update sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C4").Select()
